$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = "2025/12/04 20:00"
$ws.Range("B51").Value = "-"
$ws.Range("C51").Value = "-"
$ws.Range("D51").Value = "-"
$ws.Range("E51").Value = "-"
$ws.Range("F51").Value = "-"
$ws.Range("G51").Value = "-"
